$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Distance/MM Traveled values so they stay as text
$ws.Range("C2:C10").NumberFormat = "@"
$ws.Range("G2:G10").NumberFormat = "@"

# Row 1
$ws.Range("A1").Value = "Start"
$ws.Range("B1").Value = "End"
$ws.Range("C1").Value = "Distance"
$ws.Range("D1").Value = "Total Time"
$ws.Range("E1").Value = "Moves"
$ws.Range("F1").Value = "Time Of"
$ws.Range("G1").Value = "MM Traveled"

# Row 2
$ws.Range("A2").Value = "(314, 126)"
$ws.Range("B2").Value = "(316, 129)"
$ws.Range("C2").Value = "3.6055512754639896"
$ws.Range("D2").Value = "0:00:00.202427"
$ws.Range("E2").Value = "[]"
$ws.Range("F2").Value = "2022-07-12 14:53:32.918456"
$ws.Range("G2").Value = "2.2156694128254517"

# Row 3
$ws.Range("A3").Value = "(331, 143)"
$ws.Range("B3").Value = "(390, 144)"
$ws.Range("C3").Value = "59.008473967727724"
$ws.Range("D3").Value = "0:00:01.429046"
$ws.Range("E3").Value = "[[327, 141, datetime.timedelta(microseconds=204453), 47.523651564219136], [342, 134, datetime.timedelta(microseconds=431919), 23.550832227843234, -55.50304417350453], [378, 140, datetime.timedelta(microseconds=649382), 34.53704307117837, 16.917947900211484], [392, 150, datetime.timedelta(microseconds=845439), 12.50537915515944, -26.05943647740278], [392, 147, datetime.timedelta(seconds=1, microseconds=52929), 34.83059168124798, 421.7954717846274], [391, 145, datetime.timedelta(seconds=1, microseconds=238470), 5.762149695081819, -121.89559267902112]]"
$ws.Range("F3").Value = "2022-07-12 14:53:41.295857"
$ws.Range("G3").Value = "36.261659002748814"

# Row 4
$ws.Range("A4").Value = "(315, 138)"
$ws.Range("B4").Value = "(316, 137)"
$ws.Range("C4").Value = "1.4142135623730951"
$ws.Range("D4").Value = "0:00:00.234374"
$ws.Range("E4").Value = "[]"
$ws.Range("F4").Value = "2022-07-12 14:53:47.042532"
$ws.Range("G4").Value = "0.8690570439744343"

# Row 5
$ws.Range("A5").Value = "(312, 128)"
$ws.Range("B5").Value = "(314, 133)"
$ws.Range("C5").Value = "5.385164807134505"
$ws.Range("D5").Value = "0:00:00.222741"
$ws.Range("E5").Value = "[]"
$ws.Range("F5").Value = "2022-07-12 14:53:56.245652"
$ws.Range("G5").Value = "3.3092706314810427"

# Row 6
$ws.Range("A6").Value = "(317, 116)"
$ws.Range("B6").Value = "(269, 100)"
$ws.Range("C6").Value = "50.59644256269407"
$ws.Range("D6").Value = "0:00:02.109195"
$ws.Range("E6").Value = "[[316, 85, datetime.timedelta(microseconds=206448), 140.02749292785086], [309, 87, datetime.timedelta(microseconds=420876), 10.629603370327901, -307.4489625389021], [275, 97, datetime.timedelta(microseconds=645305), 33.749168377149104, 35.82734521942524], [268, 103, datetime.timedelta(microseconds=867770), 6.5288714421291765, -31.368100919621476], [265, 101, datetime.timedelta(seconds=1, microseconds=64247), 34.48673732353965, 435.1622002803318], [268, 103, datetime.timedelta(seconds=1, microseconds=267745), 8.275297065586479, -97.8970298528569], [271, 100, datetime.timedelta(seconds=1, microseconds=477178), 5.4637286964682, -5.8920745908618555], [271, 103, datetime.timedelta(seconds=1, microseconds=690606), 2.669464770211632, -4.046104329033586], [269, 99, datetime.timedelta(seconds=1, microseconds=905660), 3.034471739474331, 0.40302869648952055]]"
$ws.Range("F6").Value = "2022-07-12 14:54:03.428435"
$ws.Range("G6").Value = "31.092330026429746"

# Row 7
$ws.Range("A7").Value = "(289, 143)"
$ws.Range("B7").Value = "(225, 141)"
$ws.Range("C7").Value = "64.03124237432849"
$ws.Range("D7").Value = "0:00:01.275191"
$ws.Range("E7").Value = "[[261, 141, datetime.timedelta(microseconds=211903), 155.89963640477768], [254, 146, datetime.timedelta(microseconds=430346), 12.283761493726956, -333.72187707344955], [250, 146, datetime.timedelta(microseconds=631394), 3.8930755061483517, -13.289144318093939], [231, 146, datetime.timedelta(microseconds=834862), 13.985313083614901, 12.088509930343637], [225, 140, datetime.timedelta(seconds=1, microseconds=56742), 91.89563751447967, 1373.0627124680973]]"
$ws.Range("F7").Value = "2022-07-12 14:54:21.839108"
$ws.Range("G7").Value = "39.34823120099863"

# Row 8
$ws.Range("A8").Value = "(314, 132)"
$ws.Range("B8").Value = "(314, 137)"
$ws.Range("C8").Value = "5.0"
$ws.Range("D8").Value = "0:00:00.205961"
$ws.Range("E8").Value = "[]"
$ws.Range("F8").Value = "2022-07-12 14:54:26.542371"
$ws.Range("G8").Value = "3.0725806451612905"

# Row 9
$ws.Range("A9").Value = "(316, 128)"
$ws.Range("B9").Value = "(316, 135)"
$ws.Range("C9").Value = "7.0"
$ws.Range("D9").Value = "0:00:00.203456"
$ws.Range("E9").Value = "[]"
$ws.Range("F9").Value = "2022-07-12 14:54:28.700378"
$ws.Range("G9").Value = "4.301612903225807"

# Row 10
$ws.Range("A10").Value = "(314, 132)"
$ws.Range("B10").Value = "(315, 134)"
$ws.Range("C10").Value = "2.23606797749979"
$ws.Range("D10").Value = "0:00:00.205910"
$ws.Range("E10").Value = "[]"
$ws.Range("F10").Value = "2022-07-12 14:54:31.667419"
$ws.Range("G10").Value = "1.3740998377861613"

# Row 11
$ws.Range("A11").Value = "Name:"
$ws.Range("B11").Value = "m"
